$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("H17").Value = 346.1111
$ws.Range("J17").Value = 290.47058
$ws.Range("L17").Value = 871.41174
$ws.Range("N17").Value = -1207.41174
$ws.Range("H33").Value = 329.26666
$ws.Range("J33").Value = 415.6
$ws.Range("L33").Value = 415.6
$ws.Range("N33").Value = -873.6
$ws.Range("H98").Value = 2980.625
$ws.Range("I98").Value = 1977.8572
$ws.Range("K98").Value = 1977.8572
$ws.Range("M98").Value = -479.8571999999999
$ws.Range("H108").Value = 39000
$ws.Range("J108").Value = 39000
$ws.Range("L108").Value = 39000
$ws.Range("N108").Value = -46680
$ws.Range("H122").Value = 2980.625
$ws.Range("I122").Value = 1977.8572
$ws.Range("K122").Value = 5933.571599999999
$ws.Range("M122").Value = -3483.571599999999
$ws.Range("H130").Value = 19997.5
$ws.Range("J130").Value = 19997.5
$ws.Range("L130").Value = 19997.5
$ws.Range("N130").Value = -30037.5
$ws.Range("H131").Value = 11376.842
$ws.Range("I131").Value = 9120
$ws.Range("J131").Value = 16266.667
$ws.Range("K131").Value = 27360
$ws.Range("L131").Value = 48800.001
$ws.Range("M131").Value = -22320
$ws.Range("N131").Value = -58880.001
$ws.Range("H132").Value = 13884.366
$ws.Range("I132").Value = 1263.0513
$ws.Range("K132").Value = 3789.1539
$ws.Range("M132").Value = -1259.1539
$ws.Range("H138").Value = 3425.7585
$ws.Range("I138").Value = 1719.4
$ws.Range("J138").Value = 3781.25
$ws.Range("K138").Value = 5158.200000000001
$ws.Range("L138").Value = 11343.75
$ws.Range("M138").Value = -18.20000000000073
$ws.Range("N138").Value = -21623.75

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 17732.166
$ws.Range("I28").Value = 10699
$ws.Range("J28").Value = 21248.75
$ws.Range("K28").Value = 10699
$ws.Range("L28").Value = 21248.75
$ws.Range("M28").Value = -10507
$ws.Range("N28").Value = -21632.75
$ws.Range("H61").Value = 18810
$ws.Range("I61").Value = 6100
$ws.Range("K61").Value = 6100
$ws.Range("M61").Value = -5888
$ws.Range("H97").Value = 3611.2307
$ws.Range("I97").Value = 2344.5557
$ws.Range("J97").Value = 6461.25
$ws.Range("K97").Value = 2344.5557
$ws.Range("L97").Value = 6461.25
$ws.Range("M97").Value = -1848.5557
$ws.Range("N97").Value = -7453.25
$ws.Range("H99").Value = 17732.166
$ws.Range("I99").Value = 10699
$ws.Range("J99").Value = 21248.75
$ws.Range("K99").Value = 10699
$ws.Range("L99").Value = 21248.75
$ws.Range("M99").Value = -7704
$ws.Range("N99").Value = -27238.75
$ws.Range("H122").Value = 5318.077
$ws.Range("I122").Value = 5328.615
$ws.Range("K122").Value = 15985.845
$ws.Range("M122").Value = -13535.845
$ws.Range("H136").Value = 18810
$ws.Range("I136").Value = 6100
$ws.Range("K136").Value = 18300
$ws.Range("M136").Value = -15750

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H105").Value = 1555.4642
$ws.Range("I105").Value = 853.7
$ws.Range("K105").Value = 853.7
$ws.Range("M105").Value = 893.3
$ws.Range("H107").Value = 9598.714
$ws.Range("I107").Value = 9233.4375
$ws.Range("K107").Value = 9233.4375
$ws.Range("M107").Value = -7313.4375
$ws.Range("H134").Value = 4005
$ws.Range("I134").Value = 5003.6665
$ws.Range("K134").Value = 15010.9995
$ws.Range("M134").Value = -12475.9995

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 490
$ws.Range("I13").Value = 490
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 490
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -351
$ws.Range("N13").ClearContents()
$ws.Range("H31").Value = 2942.4348
$ws.Range("I31").Value = 1800.3077
$ws.Range("K31").Value = 1800.3077
$ws.Range("M31").Value = -1505.3077
$ws.Range("H34").Value = 2942.4348
$ws.Range("I34").Value = 1800.3077
$ws.Range("K34").Value = 1800.3077
$ws.Range("M34").Value = -1598.3077
$ws.Range("H94").Value = 1589.4546
$ws.Range("J94").Value = 1869.5714
$ws.Range("L94").Value = 1869.5714
$ws.Range("N94").Value = -2771.5714
$ws.Range("H134").Value = 9140.857
$ws.Range("I134").Value = 9140.857
$ws.Range("K134").Value = 27422.571
$ws.Range("M134").Value = -24887.571

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 490.64285
$ws.Range("I23").Value = 359.75
$ws.Range("K23").Value = 1079.25
$ws.Range("M23").Value = -844.25
$ws.Range("H86").Value = 452.95456
$ws.Range("I86").Value = 443.875
$ws.Range("K86").Value = 1331.625
$ws.Range("M86").Value = -145.625
$ws.Range("H89").Value = 452.95456
$ws.Range("I89").Value = 443.875
$ws.Range("K89").Value = 3994.875
$ws.Range("M89").Value = 1933.125
$ws.Range("H92").Value = 1128.4286
$ws.Range("I92").Value = 800
$ws.Range("J92").Value = 1183.1666
$ws.Range("K92").Value = 2400
$ws.Range("L92").Value = 3549.4998
$ws.Range("M92").Value = -1152
$ws.Range("N92").Value = -6045.4998

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 60
$ws.Range("J4").Value = 60
$ws.Range("L4").Value = 60
$ws.Range("N4").Value = -284
$ws.Range("H31").Value = 3174.5
$ws.Range("I31").Value = 3174.5
$ws.Range("K31").Value = 3174.5
$ws.Range("M31").Value = -2882.5
$ws.Range("H37").Value = 3174.5
$ws.Range("I37").Value = 3174.5
$ws.Range("K37").Value = 3174.5
$ws.Range("M37").Value = -2897.5
$ws.Range("H113").Value = 6675.5557
$ws.Range("I113").Value = 3046.6667
$ws.Range("J113").Value = 13933.333
$ws.Range("K113").Value = 3046.6667
$ws.Range("L113").Value = 13933.333
$ws.Range("M113").Value = -876.6667000000002
$ws.Range("N113").Value = -18273.333

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3401.75
$ws.Range("I7").Value = 2202.6667
$ws.Range("J7").Value = 6999
$ws.Range("K7").Value = 2202.6667
$ws.Range("L7").Value = 6999
$ws.Range("M7").Value = -2090.6667
$ws.Range("N7").Value = -7223
$ws.Range("H22").Value = 18182716
$ws.Range("J22").Value = 1246.5
$ws.Range("L22").Value = 1246.5
$ws.Range("N22").Value = -1836.5
$ws.Range("H27").Value = 18182716
$ws.Range("J27").Value = 1246.5
$ws.Range("L27").Value = 1246.5
$ws.Range("N27").Value = -1460.5
$ws.Range("H126").Value = 3401.75
$ws.Range("I126").Value = 2202.6667
$ws.Range("J126").Value = 6999
$ws.Range("K126").Value = 6608.000100000001
$ws.Range("L126").Value = 20997
$ws.Range("M126").Value = -4138.000100000001
$ws.Range("N126").Value = -25937
$ws.Range("H131").Value = 40000
$ws.Range("J131").Value = 40000
$ws.Range("L131").Value = 40000
$ws.Range("N131").Value = -50080
$ws.Range("H132").Value = 3086
$ws.Range("I132").Value = 2512.4
$ws.Range("K132").Value = 7537.200000000001
$ws.Range("M132").Value = -5007.200000000001

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5194.1313
$ws.Range("J132").Value = 7541.615
$ws.Range("L132").Value = 22624.845
$ws.Range("N132").Value = -27684.845
$ws.Range("H138").Value = 69500
$ws.Range("J138").Value = 69500
$ws.Range("L138").Value = 69500
$ws.Range("N138").Value = -79780
